$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (shifts former B,C to C,D)
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = 74.92

# Header for the new column
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query text (wrap text like A2)
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''German Shepherd Dog'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$ws.Range("B2").WrapText = $true

# Selection per the target worksheet XML
$ws.Range("A2").Select()
